$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1) ---
# Column E: "dependency" -> "WoS_mode"
# Column F: "WoS_state"  -> "WoS_type"
$ws.Range("E1").Value = "WoS_mode"
$ws.Range("F1").Value = "WoS_type"

# --- Update data rows 2-6: columns E (WoS_mode) and F (WoS_type) ---
$ws.Range("E2").Value = "work"
$ws.Range("F2").Value = 1

$ws.Range("E3").Value = "work"
$ws.Range("F3").Value = 1

$ws.Range("E4").Value = "study"
$ws.Range("F4").Value = 0

$ws.Range("E5").Value = "*"
$ws.Range("F5").Value = "*"

$ws.Range("E6").Value = "*"
$ws.Range("F6").Value = 1

# --- Remove old "family_types" column (G) entirely ---
$ws.Columns("G:G").Delete() | Out-Null

# --- Add new block of styled (10pt Calibri) empty cells H11:I14 ---
$ws.Range("H11:I14").Font.Size = 10

# --- Resize the now-used columns (A-H) to fit their content ---
$ws.Columns("A:H").AutoFit() | Out-Null

# --- Update selection to match the saved workbook state ---
$ws.Range("L10").Select() | Out-Null
